# Daily attendance processing - 2025-10-31 06:55:39
# Normalizes the "Recorded By" (column G) entries on the active sheet so the
# most recent recorder in the sync is listed first. Cells whose only two
# recorders are "System" and "backup@backdoor.com" are left as-is (that pair
# keeps its historical ordering); every other multi-recorder cell gets its
# first two names swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }
    if ($value -eq "") {
        continue
    }
    if ($value -eq "Recorded By") {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Length -lt 2) {
        continue
    }

    if ($parts.Length -eq 2 -and (($parts[0] -eq "System" -and $parts[1] -eq "backup@backdoor.com") -or ($parts[0] -eq "backup@backdoor.com" -and $parts[1] -eq "System"))) {
        continue
    }

    $newParts = @($parts[1], $parts[0])
    for ($i = 2; $i -lt $parts.Length; $i++) {
        $newParts += $parts[$i]
    }

    $newValue = $newParts -join ", "
    $cell.Value = $newValue
}
